$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 212.8
$ws.Range("I6").Value = 184.75
$ws.Range("K6").Value = 554.25
$ws.Range("M6").Value = -442.25

$ws.Range("H98").Value = 1054.6666
$ws.Range("I98").Value = 832.25
$ws.Range("K98").Value = 832.25
$ws.Range("M98").Value = 665.75

$ws.Range("H122").Value = 1054.6666
$ws.Range("I122").Value = 832.25
$ws.Range("K122").Value = 2496.75
$ws.Range("M122").Value = -46.75

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 19003
$ws.Range("I32").Value = 18404.8
$ws.Range("K32").Value = 18404.8
$ws.Range("M32").Value = -18117.8

$ws.Range("H45").Value = 2683.2778
$ws.Range("I45").Value = 2505.9707
$ws.Range("J45").Value = 5697.5
$ws.Range("K45").Value = 2505.9707
$ws.Range("L45").Value = 5697.5
$ws.Range("M45").Value = -2128.9707
$ws.Range("N45").Value = -6451.5

$ws.Range("H63").Value = 3424.75
$ws.Range("I63").Value = 3424.75
$ws.Range("K63").Value = 3424.75
$ws.Range("M63").Value = -2738.75

$ws.Range("H66").Value = 3424.75
$ws.Range("I66").Value = 3424.75
$ws.Range("K66").Value = 17123.75
$ws.Range("M66").Value = -13691.75

$ws.Range("H110").Value = 1262
$ws.Range("I110").Value = 1379.2
$ws.Range("J110").Value = 969
$ws.Range("K110").Value = 1379.2
$ws.Range("L110").Value = 969
$ws.Range("M110").Value = 665.8
$ws.Range("N110").Value = -5059

$ws.Range("H122").Value = 1354.8572
$ws.Range("J122").Value = 0
$ws.Range("L122").Value = 0
$ws.Range("N122").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 10495.556
$ws.Range("I107").Value = 12494.286
$ws.Range("K107").Value = 12494.286
$ws.Range("M107").Value = -10574.286

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H5").Value = 238
$ws.Range("I5").Value = 238
$ws.Range("K5").Value = 238
$ws.Range("M5").Value = -126

$ws.Range("H10").Value = 7
$ws.Range("I10").Value = 7
$ws.Range("K10").Value = 7
$ws.Range("M10").Value = 132

$ws.Range("H11").Value = 8334.166999999999
$ws.Range("I11").Value = 5
$ws.Range("J11").Value = 10000
$ws.Range("K11").Value = 5
$ws.Range("L11").Value = 10000
$ws.Range("M11").Value = 135
$ws.Range("N11").Value = -10280

$ws.Range("H31").Value = 84998.75
$ws.Range("I31").Value = 49999
$ws.Range("K31").Value = 49999
$ws.Range("M31").Value = -49704

$ws.Range("H34").Value = 84998.75
$ws.Range("I34").Value = 49999
$ws.Range("K34").Value = 49999
$ws.Range("M34").Value = -49797

$ws.Range("H107").Value = 572.1539
$ws.Range("I107").Value = 480.66666
$ws.Range("J107").Value = 650.5714
$ws.Range("K107").Value = 480.66666
$ws.Range("L107").Value = 650.5714
$ws.Range("M107").Value = 1439.33334
$ws.Range("N107").Value = -4490.5714

$ws.Range("H122").Value = 1716.5
$ws.Range("I122").Value = 1859.8
$ws.Range("J122").Value = 1000
$ws.Range("K122").Value = 5579.4
$ws.Range("L122").Value = 3000
$ws.Range("M122").Value = -3129.4
$ws.Range("N122").Value = -7900

$ws.Range("H123").Value = 0
$ws.Range("J123").Value = 0
$ws.Range("L123").Value = 0
$ws.Range("N123").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 50.714287
$ws.Range("I2").Value = 47
$ws.Range("J2").Value = 57.4
$ws.Range("K2").Value = 282
$ws.Range("L2").Value = 344.4
$ws.Range("M2").Value = -169
$ws.Range("N2").Value = -570.4

$ws.Range("H12").Value = 43.909092
$ws.Range("J12").Value = 29
$ws.Range("L12").Value = 87
$ws.Range("N12").Value = -433

$ws.Range("H32").Value = 1440.8
$ws.Range("I32").Value = 1500.5
$ws.Range("K32").Value = 4501.5
$ws.Range("M32").Value = -4218.5

$ws.Range("H33").Value = 0
$ws.Range("I33").Value = 0
$ws.Range("K33").Value = 0
$ws.Range("M33").ClearContents()

$ws.Range("H139").Value = 2073.5
$ws.Range("I139").Value = 1798.2858
$ws.Range("J139").Value = 4000
$ws.Range("K139").Value = 5394.857400000001
$ws.Range("L139").Value = 12000
$ws.Range("M139").Value = -254.8574000000008
$ws.Range("N139").Value = -22280

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 1585.7333
$ws.Range("I97").Value = 1395.6666
$ws.Range("J97").Value = 2346
$ws.Range("K97").Value = 1395.6666
$ws.Range("L97").Value = 2346
$ws.Range("M97").Value = -899.6666
$ws.Range("N97").Value = -3338

$ws.Range("H102").Value = 9499.5
$ws.Range("I102").Value = 10000
$ws.Range("J102").Value = 8999
$ws.Range("K102").Value = 10000
$ws.Range("L102").Value = 8999
$ws.Range("M102").Value = -8378
$ws.Range("N102").Value = -12243

$ws.Range("H122").Value = 10936.292
$ws.Range("I122").Value = 8211.444
$ws.Range("J122").Value = 19110.834
$ws.Range("K122").Value = 24634.332
$ws.Range("L122").Value = 57332.50199999999
$ws.Range("M122").Value = -22184.332
$ws.Range("N122").Value = -62232.50199999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 417.66666
$ws.Range("I16").Value = 417.66666
$ws.Range("K16").Value = 417.66666
$ws.Range("M16").Value = -247.66666

$ws.Range("H63").Value = 49992
$ws.Range("I63").Value = 0
$ws.Range("K63").Value = 0
$ws.Range("M63").ClearContents()

$ws.Range("H66").Value = 49992
$ws.Range("I66").Value = 0
$ws.Range("K66").Value = 0
$ws.Range("M66").ClearContents()

$ws.Range("H93").Value = 1574
$ws.Range("I93").Value = 1762.5
$ws.Range("K93").Value = 1762.5
$ws.Range("M93").Value = -514.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H53").Value = 20000
$ws.Range("J53").Value = 20000
$ws.Range("L53").Value = 20000
$ws.Range("N53").Value = -21214

$ws.Range("H122").Value = 402139
$ws.Range("I122").Value = 1000349
$ws.Range("J122").Value = 3332.3333
$ws.Range("K122").Value = 3001047
$ws.Range("L122").Value = 9996.999899999999
$ws.Range("M122").Value = -2998597
$ws.Range("N122").Value = -14896.9999

$ws.Range("H130").Value = 68990
$ws.Range("J130").Value = 68990
$ws.Range("L130").Value = 68990
$ws.Range("N130").Value = -79030

$ws.Range("H132").Value = 5621.263
$ws.Range("I132").Value = 5047.294
$ws.Range("J132").Value = 10500
$ws.Range("K132").Value = 15141.882
$ws.Range("L132").Value = 31500
$ws.Range("M132").Value = -12611.882
$ws.Range("N132").Value = -36560
